# Update of all values to match PDF edition 10 (commit 1)
# Rename the metric column headers (B1:F1) to the short/simplified
# variable-name style used in the new PDF edition.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("aircraft_operating_costs")

$ws.Range("B1").Value = "per_aircraft"
$ws.Range("C1").Value = "per_flt_hr"
$ws.Range("D1").Value = "per_flt_cycle"
$ws.Range("E1").Value = "per_seats"
$ws.Range("F1").Value = "per_ton_km"

# Match the author's final cursor position (cell selection) in the sheet.
$ws.Range("D2").Select() | Out-Null
